$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and volume-change (E) cells.
# NumberFormat is forced to text ("@") before assignment so that
# values are stored as literal strings (matching the original
# inlineStr cells) instead of being auto-converted to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.689.73'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.76%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.585.55'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -3.04%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '206.69'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -2.36%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.84%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.28'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -4.46%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.253'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -3.06%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0869'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.64%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.810.63'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -3.04%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.609.40'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.73%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.88%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.532'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -5.38%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.652.54'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.97%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -3.13%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '220.43'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -3.61%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -3.63%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.34'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -4.83%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.15'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -4.94%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.52'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -5.98%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -5.61%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.99'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.26%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.76'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.49%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.33%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.73%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.92%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.63%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0465'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.28%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -5.47%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.386.59'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.96%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -4.86%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -5.22%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.964'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -5.63%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0165'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.91%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.541'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.98%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -3.22%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.28%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.979'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.87%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.78'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -3.09%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '63.70'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -3.55%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.17'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.71%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -4.00%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.722.20'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -3.04%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '88.08'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.66%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0₆0100'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.59%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0976'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -4.82%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0500'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.87%  '
